{"js": "// The template used a Word field (fldChar begin / instrText.../ fldChar end)\n// to hold an M2Doc instruction (`m:'anydsl class diagram'.representationByName().asImage('INVALID_FORMAT')`).\n// M2Doc now expects the instruction as plain text wrapped in curly braces\n// (`{m:...}`) instead of a real Word field, so we replace the field with an\n// equivalent run of literal text runs (preserving the original run\n// boundaries / formatting) reading \"{m:...asImage('INVALID_FORMAT')}\".\n\nconst paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph/field holding the M2Doc \"asImage\" instruction.\nlet targetPara = null;\nlet targetField = null;\nfor (const p of paras.items) {\n  const fields = p.fields;\n  fields.load(\"items\");\n  await context.sync();\n  for (const f of fields.items) {\n    f.load(\"code\");\n    await context.sync();\n    if (f.code && f.code.indexOf(\"asImage\") !== -1) {\n      targetPara = p;\n      targetField = f;\n      break;\n    }\n  }\n  if (targetField) break;\n}\n\nif (!targetField) {\n  throw new Error(\"Could not find the asImage(...) field in the document.\");\n}\n\n// The literal text chunks that used to be separate <w:instrText> runs inside\n// the field; reproduced here as plain text runs, with the instruction now\n// fenced in \"{ ... }\" braces as a literal run of text.\nconst textChunks = [\n  \"{m:\",\n  \"'\",\n  \"anydsl class diagram\",\n  \"'.\",\n  \"r\",\n  \"epresentation\",\n  \"By\",\n  \"Name\",\n  \"()\",\n  \".\",\n  \"asImage('INVALID_FORMAT')}\",\n];\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\nconst runsXml = textChunks\n  .map(\n    (t) =>\n      `<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>${escapeXml(\n        t\n      )}</w:t></w:r>`\n  )\n  .join(\"\");\n\nconst ooxmlFrag =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p>\" +\n  runsXml +\n  \"</w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\n// Insert the literal-text runs right before the field, then delete the\n// field itself (which removes the fldChar begin/end + instrText runs).\ntargetPara.insertOoxml(ooxmlFrag, Word.InsertLocation.start);\ntargetField.delete();\nawait context.sync();\n", "ps1": "# The template used a real Word field (fldChar begin / instrText... / fldChar\n# end) to hold an M2Doc instruction\n# (m:'anydsl class diagram'.representationByName().asImage('INVALID_FORMAT')).\n# M2Doc now expects the instruction as plain text wrapped in curly braces\n# (\"{m:...}\") instead of a real Word field, so we replace the field with an\n# equivalent run of literal text runs (preserving the original run\n# boundaries / character formatting) reading\n# \"{m:...asImage('INVALID_FORMAT')}\".\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph/field holding the M2Doc \"asImage\" instruction.\n$targetPara = $null\n$targetField = $null\nforeach ($p in $d.Paragraphs) {\n    foreach ($fld in $p.Range.Fields) {\n        if ($fld.Code.Text -like \"*asImage*\") {\n            $targetPara = $p\n            $targetField = $fld\n            break\n        }\n    }\n    if ($targetField) { break }\n}\n\nif (-not $targetField) {\n    throw \"Could not find the asImage(...) field in the document.\"\n}\n\n# The literal text chunks that used to be separate <w:instrText> runs inside\n# the field; reproduced here as plain text runs, with the instruction now\n# fenced in \"{ ... }\" braces as a literal run of text.\n$textChunks = @(\n    \"{m:\",\n    \"'\",\n    \"anydsl class diagram\",\n    \"'.\",\n    \"r\",\n    \"epresentation\",\n    \"By\",\n    \"Name\",\n    \"()\",\n    \".\",\n    \"asImage('INVALID_FORMAT')}\"\n)\n\n$runsXml = \"\"\nforeach ($chunk in $textChunks) {\n    $escaped = $chunk -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'\n    $runsXml += \"<w:r><w:rPr><w:sz w:val=`\"24`\"/><w:szCs w:val=`\"24`\"/></w:rPr><w:t>$escaped</w:t></w:r>\"\n}\n\n$ooxmlFrag = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' + $runsXml + '</w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n# Insert the literal-text runs right before the field (a collapsed range at\n# the paragraph start keeps the paragraph's own <w:pPr> / paraId intact),\n# then delete the field itself (which removes the fldChar begin/end +\n# instrText runs).\n$insertionPoint = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)\n$insertionPoint.InsertXML($ooxmlFrag)\n$targetField.Delete()\n"}
